$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 44 (rows 45 and 46), pushing the
# former rows 45..74 down to 47..76. This mirrors a new weekly batch of
# price records being added to the dataset.
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

# Row 45 - new record
$ws.Range("A45").Value = 9
$ws.Range("B45").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C45").Value = "Metropolitana"
$ws.Range("D45").Value = 44586
$ws.Range("E45").Value = 13
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100101
$ws.Range("H45").Value = "Berries"
$ws.Range("I45").Value = 100101004
$ws.Range("J45").Value = "Frambuesa"
$ws.Range("K45").Value = "Sin especificar"
$ws.Range("L45").Value = "Especial"
$ws.Range("M45").Value = 250
$ws.Range("N45").Value = 8000
$ws.Range("O45").Value = 8000
$ws.Range("P45").Value = 8000
$ws.Range("Q45").Value = "$/bandeja 2 kilos"
$ws.Range("R45").Value = "Provincia de Linares"
$ws.Range("S45").Value = 4000
$ws.Range("T45").Value = 2

# Row 46 - new record
$ws.Range("A46").Value = 9
$ws.Range("B46").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C46").Value = "Metropolitana"
$ws.Range("D46").Value = 44586
$ws.Range("E46").Value = 13
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100101
$ws.Range("H46").Value = "Berries"
$ws.Range("I46").Value = 100101004
$ws.Range("J46").Value = "Frambuesa"
$ws.Range("K46").Value = "Sin especificar"
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 200
$ws.Range("N46").Value = 7000
$ws.Range("O46").Value = 7000
$ws.Range("P46").Value = 7000
$ws.Range("Q46").Value = "$/bandeja 2 kilos"
$ws.Range("R46").Value = "Provincia de Linares"
$ws.Range("S46").Value = 3500
$ws.Range("T46").Value = 2
